$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51
$ws.Cells.Item($row, 1).Value = "Daniele Ruzzenenti"
$ws.Cells.Item($row, 2).Value = "Elia Battisti | U.SGUARNA"
$ws.Cells.Item($row, 3).Value = "Michele Merighi | Clitoriders"
$ws.Cells.Item($row, 4).Value = "Giacomo  Gasparini  | Mai una gioia"
$ws.Cells.Item($row, 5).Value = "Manuel Emanuelli | SdrumALA"
$ws.Cells.Item($row, 6).Value = "Moris Benedetti | Gli Introvabili"
